$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.292
$ws.Range("C7").Value = -12.909
$ws.Range("E7").Value = 16.374
$ws.Range("E15").Value = 16.21
$ws.Range("C16").Value = -13.377
$ws.Range("E21").Value = 16.557
$ws.Range("E22").Value = 16.647
$ws.Range("E23").Value = 16.434
$ws.Range("C28").Value = -13.219
$ws.Range("C29").Value = -11.736
$ws.Range("C32").Value = -13.434
$ws.Range("E34").Value = 16.907
$ws.Range("C40").Value = -12.232
$ws.Range("E43").Value = 17.219
$ws.Range("E45").Value = 16.696
$ws.Range("E50").Value = 16.433
$ws.Range("E51").Value = 16.77
$ws.Range("C52").Value = -11.355
$ws.Range("C57").Value = -13.829
$ws.Range("C66").Value = -11.598
$ws.Range("E66").Value = 17.157
$ws.Range("E67").Value = 17.398
$ws.Range("E79").Value = 17.017
$ws.Range("E84").Value = 16.492
$ws.Range("E92").Value = 17.797
$ws.Range("E97").Value = 16.77
$ws.Range("C100").Value = -13.108
